# Updates the Enterprise Java time log through 2/17, and collapses the
# trailing four "scratch" bullet rows (22-25) into a single consolidated
# entry that now lives in the new row 21 (with a real date/hours pair),
# leaving row 22 as a blank spacer row beneath it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 21: date, hours, and task note for 2/17 -----------------------
$ws.Cells.Item(21, 1).NumberFormat = $ws.Cells.Item(20, 1).NumberFormat
$ws.Cells.Item(21, 1).Value = Get-Date -Year 2019 -Month 2 -Day 17 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(21, 2).Value = 6.5
$ws.Cells.Item(21, 4).Value = "Indie Project: Revised user stories and problem statement; finished screen designs; added work to GitHub and appropriate links to README.  Week 4: Watched intro video.  Professional Development: signed up to present on Hibernate Search; also, will this help users search across forums or perhaps across site?  Would it make admin related tasks easier?"
$ws.Rows.Item(21).RowHeight = 45

# --- Old rows 22-25 collapse down to a single blank spacer row 22 ----------
$ws.Cells.Item(22, 4).Clear()
$ws.Rows("23:25").Delete()
$ws.Cells.Item(22, 1).NumberFormat = $ws.Cells.Item(21, 1).NumberFormat

# --- Selection / view state, matching the saved workbook -------------------
$ws.Range("D24:D29").Select()
